$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel keeps them as text instead of silently converting to a number.
$numericLookingCells = @("D5","D6","D10","D11","D17","D19","D20","D22","D29","D31","D36","D39","D40","D41","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '61.927.40'
$ws.Range("E2").Value = '  +3.20%  '
$ws.Range("D3").Value = '3.415.13'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '577.56'
$ws.Range("E5").Value = '  +2.64%  '
$ws.Range("D6").Value = '138.98'
$ws.Range("E6").Value = '  +8.27%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.412.02'
$ws.Range("E8").Value = '  +3.48%  '
$ws.Range("E9").Value = '  +1.73%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.128'
$ws.Range("E10").Value = '  +10.67%  '
$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '7.49'
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("E12").Value = '  +7.12%  '
$ws.Range("D13").Value = '4.000.65'
$ws.Range("E13").Value = '  +3.76%  '
$ws.Range("E14").Value = '  +2.11%  '
$ws.Range("E15").Value = '  +9.52%  '
$ws.Range("D16").Value = '3.422.61'
$ws.Range("E16").Value = '  +3.84%  '
$ws.Range("D17").Value = '25.53'
$ws.Range("E17").Value = '  +6.31%  '
$ws.Range("D18").Value = '61.985.54'
$ws.Range("E18").Value = '  +2.96%  '
$ws.Range("D19").Value = '14.15'
$ws.Range("E19").Value = '  +7.06%  '
$ws.Range("D20").Value = '5.92'
$ws.Range("E20").Value = '  +5.88%  '
$ws.Range("E21").Value = '  +7.36%  '
$ws.Range("D22").Value = '391.10'
$ws.Range("E22").Value = '  +12.31%  '
$ws.Range("E23").Value = '  +4.26%  '
$ws.Range("D24").Value = '3.553.00'
$ws.Range("E24").Value = '  +3.77%  '
$ws.Range("E25").Value = '  +20.06%  '
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("E27").Value = '  +4.47%  '
$ws.Range("E28").Value = '  +10.71%  '
$ws.Range("D29").Value = '7.66'
$ws.Range("E29").Value = '  +5.36%  '
$ws.Range("E30").Value = '  +0.07%  '
$ws.Range("D31").Value = '8.32'
$ws.Range("E31").Value = '  +6.81%  '
$ws.Range("E32").Value = '  +6.15%  '
$ws.Range("E33").Value = '  +4.00%  '
$ws.Range("D34").Value = '3.448.80'
$ws.Range("E34").Value = '  +3.77%  '
$ws.Range("D36").Value = '23.62'
$ws.Range("E36").Value = '  +4.64%  '
$ws.Range("E37").Value = '  +4.15%  '
$ws.Range("E38").Value = '  +4.16%  '
$ws.Range("D39").Value = '1.58'
$ws.Range("E39").Value = '  +7.44%  '
$ws.Range("D40").Value = '162.30'
$ws.Range("E40").Value = '  +3.74%  '
$ws.Range("D41").Value = '0.0794'
$ws.Range("E41").Value = '  +6.37%  '
$ws.Range("E42").Value = '  +15.76%  '
$ws.Range("E43").Value = '  +7.52%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").Value = '25.47'
$ws.Range("E44").Value = '  +13.36%  '
$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '1.23'
$ws.Range("E46").Value = '  +6.25%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '4.49'
$ws.Range("E47").Value = '  +4.66%  '
$ws.Range("D48").Value = '41.68'
$ws.Range("E48").Value = '  +2.64%  '
$ws.Range("D49").Value = '6.98'
$ws.Range("E49").Value = '  +4.49%  '
$ws.Range("D50").Value = '23.24'
$ws.Range("E50").Value = '  +6.99%  '
$ws.Range("D51").Value = '2.402.27'
$ws.Range("E51").Value = '  +12.00%  '
